$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert the two new rows in the middle of the table first (structure only) ---
# New row for "Simple if else example", right after row 63.
$ws.Rows.Item(64).Insert()
# New row for "Find the number of days...", right after the (now shifted) row 65.
$ws.Rows.Item(66).Insert()

# --- Fill in the new/changed cells in the same order the original author did ---
$ws.Range("B66").Value = "int(timestamp(now) - timestamp(user.meta.created)) / 86400"
$ws.Range("A66").Value = "Find the number of days the user account was created`nThe function int() is used to convert the duration type to a integer"
$ws.Range("B66").Style = "Normal"
$ws.Rows.Item(66).RowHeight = 28.8

$ws.Range("A64").Value = "Simple if else example"
$ws.Range("B64").Value = "statements:`n    - if:`n        match: false`n        block:`n            - return: string(`"hello`")`n        else:`n            - return: string(`"goodbye`")"
$ws.Rows.Item(64).RowHeight = 100.8

# --- Update the title cell (A1, merged A1:B1) ---
# Keep "Sample Scenarios" bold/italic/underline (inherited from the cell style)
# and the description line in regular weight, now with extra text appended.
$line1 = "Sample Scenarios`n"
$line2 = "This file contains information on how to reference the user object to get user information using advanced attributes as well as simple examples"
$titleCell = $ws.Range("A1")
$titleCell.Value = $line1 + $line2
$descRun = $titleCell.Characters($line1.Length + 1, $line2.Length)
$descRun.Font.Bold = $false
$descRun.Font.Italic = $false
$descRun.Font.Underline = $false

# --- Append two new rows at the end for the "Simple multi-line rule" scenarios ---
$ws.Range("A67").Value = "Simple multi-line rule to return a string value"
$ws.Range("B67").Value = "statements:`n    - return: string(`"some value`")"
$ws.Range("A67:B67").WrapText = $true
$ws.Rows.Item(67).RowHeight = 28.8

$ws.Range("A68").Value = "Simple multi-line rule to return a boolean value"
$ws.Range("B68").Value = "statements:`n    - return: true"
$ws.Range("A68:B68").WrapText = $true
$ws.Rows.Item(68).RowHeight = 28.8

# --- Update the view state to match the new bottom of the sheet ---
$ws.Range("B70").Select() | Out-Null
